# Insert a new row above the current row 149. Excel will shift rows
# 149..186 down to 150..187 automatically (carrying values & styles),
# which matches the target state of the workbook (each record moved
# down by one row, last old record becoming new row 187).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("149:149").Insert()

# Populate the newly inserted row 149 with its data.
$ws.Range("A149").Value = 3
$ws.Range("B149").Value = "Femacal de La Calera"
$ws.Range("C149").Value = "Coquimbo"
$ws.Range("D149").Value = 44798
$ws.Range("E149").Value = 5
$ws.Range("F149").Value = 100112026
$ws.Range("G149").Value = "Haba"
$ws.Range("H149").Value = "Sin especificar"
$ws.Range("I149").Value = "Primera"
$ws.Range("J149").Value = 125
$ws.Range("K149").Value = 14000
$ws.Range("L149").Value = 15000
$ws.Range("M149").Value = 14480
$ws.Range("N149").Value = "`$/malla 25 kilos"
$ws.Range("O149").Value = "Provincia de Limarí"
$ws.Range("P149").Value = 579
$ws.Range("Q149").Value = 25
$ws.Range("R149").Value = "Hortaliza"
